# "Working on RT quantile of Model Preds"
# Appends a new pair of fitted-model rows (7 & 8) to both the "Mix" and "VP"
# sheets, restyles the new Mix!M8 cell with the workbook's existing
# scientific-notation style, and moves the active tab/selection so that
# Mix!M8 is the active cell (VP keeps its own remembered A8:B8 selection).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Mix
$ws2 = $wb.Worksheets.Item(2)   # VP

# --- Mix sheet: new rows 7 and 8 (A:N) ---
$ws1.Cells.Item(7,1).Value = 25407.294767784
$ws1.Cells.Item(7,2).Value = 50921.2430247383
$ws1.Cells.Item(7,3).Value = 2.84815427980301
$ws1.Cells.Item(7,4).Value = -0.00214561308616434
$ws1.Cells.Item(7,5).Value = 3.04494116111179
$ws1.Cells.Item(7,6).Value = -0.000593477090334862
$ws1.Cells.Item(7,7).Value = 0.0137083179388561
$ws1.Cells.Item(7,8).Value = 0.0175233303509455
$ws1.Cells.Item(7,9).Value = 2.58279388315756
$ws1.Cells.Item(7,10).Value = 1.37636969257252
$ws1.Cells.Item(7,11).Value = 0.53153711445798
$ws1.Cells.Item(7,12).Value = 0.559579074820132
$ws1.Cells.Item(7,13).Value = 0.0167000024371789
$ws1.Cells.Item(7,14).Value = 1.18336697194636

$ws1.Cells.Item(8,1).Value = 12704.991302479
$ws1.Cells.Item(8,2).Value = 25515.6200897279
$ws1.Cells.Item(8,3).Value = 2.6604749466414
$ws1.Cells.Item(8,4).Value = -0.00278331185857538
$ws1.Cells.Item(8,5).Value = 2.63550043743625
$ws1.Cells.Item(8,6).Value = 0.0123075694016193
$ws1.Cells.Item(8,7).Value = 0.00725041247253861
$ws1.Cells.Item(8,8).Value = 0.00241637765593459
$ws1.Cells.Item(8,9).Value = 2.45956850608643
$ws1.Cells.Item(8,10).Value = 1.38438596728826
$ws1.Cells.Item(8,11).Value = 0.575306567223471
$ws1.Cells.Item(8,12).Value = 0.618873956879957
$ws1.Cells.Item(8,13).Value = 0.0000274082775257352
$ws1.Cells.Item(8,14).Value = 1.21161371469806

# M8 uses the scientific-notation style already present on H2:H4 (style index 1)
$ws1.Cells.Item(8,13).NumberFormat = "0.00E+00"

# --- VP sheet: new rows 7 and 8 (A:K) ---
$ws2.Cells.Item(7,1).Value = 23810.33483076
$ws2.Cells.Item(7,2).Value = 47700.6597783977
$ws2.Cells.Item(7,3).Value = 1.85819385195165
$ws2.Cells.Item(7,4).Value = 0.0109097621350335
$ws2.Cells.Item(7,5).Value = 1.56313097812632
$ws2.Cells.Item(7,6).Value = 0.692899797565077
$ws2.Cells.Item(7,7).Value = 1.129987814522
$ws2.Cells.Item(7,8).Value = 0.00208615711603766
$ws2.Cells.Item(7,9).Value = 1.46244802337094
$ws2.Cells.Item(7,10).Value = 0.624099110078687
$ws2.Cells.Item(7,11).Value = 1.87993837822989

$ws2.Cells.Item(8,1).Value = 13784.9478595788
$ws2.Cells.Item(8,2).Value = 27649.123832735
$ws2.Cells.Item(8,3).Value = 1.85819385195165
$ws2.Cells.Item(8,4).Value = 0.0109097621350335
$ws2.Cells.Item(8,5).Value = 1.56313097812632
$ws2.Cells.Item(8,6).Value = 0.692899797565077
$ws2.Cells.Item(8,7).Value = 1.129987814522
$ws2.Cells.Item(8,8).Value = 0.00208615711603766
$ws2.Cells.Item(8,9).Value = 1.46244802337094
$ws2.Cells.Item(8,10).Value = 0.624099110078687
$ws2.Cells.Item(8,11).Value = 1.87993837822989

# --- Selection / active-sheet bookkeeping ---
# Remember VP's selection (A8:B8) before switching the active tab to Mix,
# so VP's sheetView keeps that selection without being the active tab.
$ws2.Range("A8:B8").Select()
$ws1.Range("M8").Select()
